$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '43.971.97'
$ws.Range("E2").Value = '  +2.56%  '
# Row 3
Set-TextValue "D3" '2.341.60'
$ws.Range("E3").Value = '  +2.57%  '
# Row 4
$ws.Range("E4").Value = '  -0.18%  '
# Row 5
Set-TextValue "D5" '312.80'
$ws.Range("E5").Value = '  -0.11%  '
# Row 6
Set-TextValue "D6" '108.37'
$ws.Range("E6").Value = '  +3.17%  '
# Row 7
$ws.Range("E7").Value = '  +1.24%  '
# Row 8
$ws.Range("E8").Value = '  -0.15%  '
# Row 9
$ws.Range("E9").Value = '  +3.13%  '
# Row 10
Set-TextValue "D10" '41.27'
$ws.Range("E10").Value = '  +4.27%  '
# Row 11
Set-TextValue "D11" '0.0919'
$ws.Range("E11").Value = '  +2.04%  '
# Row 12
Set-TextValue "D12" '8.58'
$ws.Range("E12").Value = '  +2.55%  '
# Row 13
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D13" '1.01'
$ws.Range("E13").Value = '  +2.17%  '
# Row 14
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue "D14" '0.108'
$ws.Range("E14").Value = '  -0.99%  '
# Row 15
Set-TextValue "D15" '15.51'
$ws.Range("E15").Value = '  +2.23%  '
# Row 16
Set-TextValue "D16" '2.697.32'
$ws.Range("E16").Value = '  +2.52%  '
# Row 17
Set-TextValue "D17" '2.342.60'
$ws.Range("E17").Value = '  +2.75%  '
# Row 18
Set-TextValue "D18" '43.917.10'
$ws.Range("E18").Value = '  +2.67%  '
# Row 19
Set-TextValue "D19" '7.56'
$ws.Range("E19").Value = '  +2.82%  '
# Row 20
$ws.Range("E20").Value = '  +1.86%  '
# Row 21
$ws.Range("E21").Value = '  -6.27%  '
# Row 22
Set-TextValue "D22" '74.29'
$ws.Range("E22").Value = '  +0.78%  '
# Row 23
$ws.Range("E23").Value = '  -1.29%  '
# Row 24
Set-TextValue "D24" '268.79'
$ws.Range("E24").Value = '  +1.73%  '
# Row 25
Set-TextValue "D25" '2.29'
$ws.Range("E25").Value = '  +3.86%  '
# Row 26
$ws.Range("E26").Value = '  -0.12%  '
# Row 27
$ws.Range("E27").Value = '  +5.45%  '
# Row 28
$ws.Range("E28").Value = '  +3.37%  '
# Row 29
$ws.Range("E29").Value = '  -1.89%  '
# Row 30
Set-TextValue "D30" '39.03'
$ws.Range("E30").Value = '  +5.82%  '
# Row 31
Set-TextValue "D31" '22.66'
$ws.Range("E31").Value = '  +0.98%  '
# Row 32
Set-TextValue "D32" '168.52'
$ws.Range("E32").Value = '  +1.14%  '
# Row 33
Set-TextValue "D33" '0.0892'
$ws.Range("E33").Value = '  +2.64%  '
# Row 34
Set-TextValue "D34" '2.81'
$ws.Range("E34").Value = '  +8.64%  '
# Row 35
$ws.Range("E35").Value = '  +1.57%  '
# Row 36
Set-TextValue "D36" '4.76'
$ws.Range("E36").Value = '  +5.28%  '
# Row 37
Set-TextValue "D37" '0.114'
$ws.Range("E37").Value = '  +0.94%  '
# Row 38
Set-TextValue "D38" '0.0364'
$ws.Range("E38").Value = '  +4.25%  '
# Row 39
Set-TextValue "D39" '2.91'
$ws.Range("E39").Value = '  +9.71%  '
# Row 40
$ws.Range("E40").Value = '  +1.54%  '
# Row 42
Set-TextValue "D42" '104.51'
$ws.Range("E42").Value = '  +10.00%  '
# Row 43
Set-TextValue "D43" '0.240'
$ws.Range("E43").Value = '  +4.07%  '
# Row 44
Set-TextValue "D44" '71.65'
$ws.Range("E44").Value = '  +2.16%  '
# Row 45
Set-TextValue "D45" '13.31'
$ws.Range("E45").Value = '  +9.83%  '
# Row 46
$ws.Range("E46").Value = '  +0.04%  '
# Row 47
Set-TextValue "D47" '114.46'
$ws.Range("E47").Value = '  +1.08%  '
# Row 48
Set-TextValue "D48" '1.666.20'
$ws.Range("E48").Value = '  -3.68%  '
# Row 49
Set-TextValue "D49" '76.92'
$ws.Range("E49").Value = '  -2.17%  '
# Row 50
Set-TextValue "D50" '8.98'
$ws.Range("E50").Value = '  +3.04%  '
# Row 51
$ws.Range("E51").Value = '  +8.87%  '
